# Write a small "Number"/"Data" list out to the worksheet (twice, as the
# source list was appended to the sheet in two passes) and leave the
# workbook positioned/selected the way the authoring run left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")
$ws.Activate()

# The list of records that got written to the sheet.
$records = @()
for ($i = 0; $i -lt 10; $i++) {
    $records += ,@($i, ("{0:D2}DataAA" -f $i))
}

$row = 2
for ($pass = 0; $pass -lt 2; $pass++) {
    $ws.Cells.Item($row, 1).Value = "Number"
    $ws.Cells.Item($row, 2).Value = "Data"
    $row++

    foreach ($rec in $records) {
        $ws.Cells.Item($row, 1).Value = $rec[0]
        $ws.Cells.Item($row, 2).Value = $rec[1]
        $row++
    }
}

$ws.Range("J16").Select() | Out-Null
